$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# The workbook currently has two tabs: "Sheet2" (contains the real BOM
# content + the query table) and "Sheet1" (a mysterious blank sheet).
# Grab the data sheet first (it's the active one), delete the blank
# sheet, and rename the data sheet to "Sheet1".
$ws = $wb.ActiveSheet
$wb.Worksheets("Sheet1").Delete()
$ws.Name = "Sheet1"

# Add "Date" / "Notes" headers next to the VERSION INFO block, and record
# the first BOM revision's date + a note.
$ws.Range("B19").Value = "Date"
$ws.Range("C19").Value = "Notes"

$ws.Range("B21").Value = (Get-Date -Year 2020 -Month 8 -Day 28 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C21").Value = "First version of BOM"

# Reset the view: scroll back to the top and move the selection onto the
# newly-added Notes cell.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F19").Select()
